# Apply "Add Batteries in addition to Ultracaps" change
# to the gStation worksheet of the eco_cost_inputs workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("gStation")

# --- Change drivetrain_type value (row 9) ---
$ws.Range("B9").Value = 1

# --- Insert a new row for "elecSto_type" right before the old
#     "ultracap.p" row (row 20) ---
$ws.Rows.Item(20).Insert()
$ws.Range("A20").Value = "elecSto_type"
$ws.Range("B20").Value = 1

# --- Insert two new rows after "ultracap.N" (now row 22) for the
#     new battery parameters "batt.p" and "batt.N" ---
$ws.Rows.Item(23).Insert()
$ws.Rows.Item(23).Insert()

$ws.Range("A23").Value = "batt.p"
$ws.Range("B23").Value = 180
$ws.Range("B23").NumberFormat = "0.00E+00"

$ws.Range("A24").Value = "batt.N"
$ws.Range("B24").Formula = "=10000"

# --- Update sheet view: activating "gStation" makes it the
#     selected tab (and automatically clears the previous
#     tabSelected flag on "tether"), then move the selection
#     to F20 ---
$ws.Activate()
$ws.Range("F20").Select()
